$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (data, loja) before the existing "nome" column,
# shifting all existing columns A:G to C:I.
$ws.Columns("A:B").Insert()

# The newly inserted columns have no formatting; copy the header formatting
# from the (now shifted) "nome" header cell C1 so A1:B1 match the other
# header cells (bold, centered, top-aligned, thin border).
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("A1").Value = "data"
$ws.Range("B1").Value = "loja"

$ws.Range("A2").Value = "30/07/2024"
$ws.Range("B2").Value = "jotacar parts"
$ws.Range("C2").Value = "Fonte Carregador Jfa 70a Bivolt Com Medidor Cca"
$ws.Range("D2").Value = "FONTE 70A STORM"
$ws.Range("F2").Value = "Acima"
$ws.Range("G2").Value = "FULL"
$ws.Range("H2").Value = "classico"
$ws.Range("I2").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-70a-bivolt-com-medidor-cca/p/MLB21455208?pdp_filters=seller_id:129467482#searchVariation=MLB21455208&position=5&search_layout=stack&type=product&tracking_id=f4dc8cf6-1ac8-42d8-ad0d-bcd4f0ed286c"
$ws.Range("E2").Value = 499.9
$ws.Range("A3").Value = "30/07/2024"
$ws.Range("B3").Value = "jotacar parts"
$ws.Range("C3").Value = "Fonte Automotiva 120a Amperes Jfa Carregador Cor Preto"
$ws.Range("D3").Value = "FONTE 120A STORM"
$ws.Range("F3").Value = "Acima"
$ws.Range("G3").Value = "NA"
$ws.Range("H3").Value = "premium"
$ws.Range("I3").Value = "https://www.mercadolivre.com.br/fonte-automotiva-120a-amperes-jfa-carregador-cor-preto/p/MLB21392652?pdp_filters=seller_id:129467482#searchVariation=MLB21392652&position=3&search_layout=stack&type=product&tracking_id=f4dc8cf6-1ac8-42d8-ad0d-bcd4f0ed286c"
$ws.Range("E3").Value = 674.97
$ws.Range("A4").Value = "30/07/2024"
$ws.Range("B4").Value = "jotacar parts"
$ws.Range("C4").Value = "Fonte Carregador Automotivo Jfa 120a Sci Bivolt Pwm A Melhor Cor Preto"
$ws.Range("D4").Value = "FONTE 120A STORM"
$ws.Range("F4").Value = "Acima"
$ws.Range("G4").Value = "FULL"
$ws.Range("H4").Value = "classico"
$ws.Range("I4").Value = "https://www.mercadolivre.com.br/fonte-carregador-automotivo-jfa-120a-sci-bivolt-pwm-a-melhor-cor-preto/p/MLB27869459?pdp_filters=seller_id:129467482#searchVariation=MLB27869459&position=4&search_layout=stack&type=product&tracking_id=f4dc8cf6-1ac8-42d8-ad0d-bcd4f0ed286c"
$ws.Range("E4").Value = 634.6
$ws.Range("A5").Value = "30/07/2024"
$ws.Range("B5").Value = "jotacar parts"
$ws.Range("C5").Value = "Fonte Carregador Jfa 60a Bivolt Storm Com Medidor Cca"
$ws.Range("D5").Value = "FONTE 60A STORM"
$ws.Range("F5").Value = "Acima"
$ws.Range("G5").Value = "NA"
$ws.Range("H5").Value = "premium"
$ws.Range("I5").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-60a-bivolt-storm-com-medidor-cca/p/MLB21320712?pdp_filters=seller_id:129467482#searchVariation=MLB21320712&position=1&search_layout=stack&type=product&tracking_id=f4dc8cf6-1ac8-42d8-ad0d-bcd4f0ed286c"
$ws.Range("E5").Value = 489.9
$ws.Range("A6").Value = "30/07/2024"
$ws.Range("B6").Value = "jotacar parts"
$ws.Range("C6").Value = "Fonte Carregador De Bateria 70a Storm Plus Automotiva Jfa"
$ws.Range("D6").Value = "FONTE 70A STORM"
$ws.Range("F6").Value = "Acima"
$ws.Range("G6").Value = "NA"
$ws.Range("H6").Value = "premium"
$ws.Range("I6").Value = "https://produto.mercadolivre.com.br/MLB-3703963247-fonte-carregador-de-bateria-70a-storm-plus-automotiva-jfa-_JM#position%3D6%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df4dc8cf6-1ac8-42d8-ad0d-bcd4f0ed286c"
$ws.Range("E6").Value = 539.9
$ws.Range("A7").Value = "30/07/2024"
$ws.Range("B7").Value = "jotacar parts"
$ws.Range("C7").Value = "Fonte Bivolt Para Caixa Bob 90a Jfa Carregador De Bateria"
$ws.Range("D7").Value = "FONTE 90 BOB"
$ws.Range("F7").Value = "Acima"
$ws.Range("G7").Value = "NA"
$ws.Range("H7").Value = "premium"
$ws.Range("I7").Value = "https://produto.mercadolivre.com.br/MLB-3703784627-fonte-bivolt-para-caixa-bob-90a-jfa-carregador-de-bateria-_JM#position%3D8%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df4dc8cf6-1ac8-42d8-ad0d-bcd4f0ed286c"
$ws.Range("E7").Value = 446.3
